$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A value to be stored as text (like "2025-02-28" in the rows
# above) instead of being auto-converted to an Excel date serial number.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "2025-03-18"
# Reset the cell style back to Normal so no stray number-format styling is
# left on the new cell (matches the un-styled cells used by the other rows).
$ws.Range("A22").Style = "Normal"

$ws.Range("B22").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C22").Value = "NA"
$ws.Range("D22").Value = 1
